# Generate Report for Handback
# Updates the "Correspond Handoff Datetime" (column E) and
# "Correspond Handback DateTime" (column H) values for the first data
# row (the 3a249cac... file) on both the "zh-cn" and "de-de" sheets,
# reflecting a newer handoff/handback cycle that was just generated.

$wb = $excel.ActiveWorkbook

$ws_zhcn = $wb.Worksheets.Item("zh-cn")
$ws_zhcn.Range("E2").Value = "2016-03-20 16:42:56"
$ws_zhcn.Range("H2").Value = "2016-03-20 16:43:22"

$ws_dede = $wb.Worksheets.Item("de-de")
$ws_dede.Range("E2").Value = "2016-03-20 16:42:59"
$ws_dede.Range("H2").Value = "2016-03-20 16:43:28"
